# repull data, push all data, mean calculation
# Update the dSF (column F) values to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    10 = 0
    13 = -7
    14 = -3
    17 = -2
    19 = -5
    20 = -6
    22 = 5
    23 = -5
    25 = 2
    26 = 1
    28 = -1
    30 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
